# Auto-applied numeric updates to Leve profit/price tables, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2546.4482
$ws.Range("I62").Value = 2024.7894
$ws.Range("J62").Value = 3537.6
$ws.Range("K62").Value = 2024.7894
$ws.Range("L62").Value = 3537.6
$ws.Range("M62").Value = -1400.7894
$ws.Range("N62").Value = -4785.6
$ws.Range("H65").Value = 2546.4482
$ws.Range("I65").Value = 2024.7894
$ws.Range("J65").Value = 3537.6
$ws.Range("K65").Value = 10123.947
$ws.Range("L65").Value = 17688
$ws.Range("M65").Value = -7003.947
$ws.Range("N65").Value = -23928
$ws.Range("H132").Value = 1383.6666
$ws.Range("I132").Value = 1383.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4150.9998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1620.9998
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 577308.7
$ws.Range("I137").Value = 3153.3635
$ws.Range("J137").Value = 928181.4399999999
$ws.Range("K137").Value = 9460.0905
$ws.Range("L137").Value = 2784544.32
$ws.Range("M137").Value = -6910.0905
$ws.Range("N137").Value = -2789644.32
$ws.Range("H138").Value = 5732.4863
$ws.Range("I138").Value = 3070.7
$ws.Range("J138").Value = 6718.3335
$ws.Range("K138").Value = 9212.099999999999
$ws.Range("L138").Value = 20155.0005
$ws.Range("M138").Value = -4072.099999999999
$ws.Range("N138").Value = -30435.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13412.25
$ws.Range("I32").Value = 13926.072
$ws.Range("K32").Value = 13926.072
$ws.Range("M32").Value = -13639.072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1783.4509
$ws.Range("I86").Value = 1620.1333
$ws.Range("J86").Value = 3008.3333
$ws.Range("K86").Value = 1620.1333
$ws.Range("L86").Value = 3008.3333
$ws.Range("M86").Value = -497.1333
$ws.Range("N86").Value = -5254.3333
$ws.Range("H89").Value = 1783.4509
$ws.Range("I89").Value = 1620.1333
$ws.Range("J89").Value = 3008.3333
$ws.Range("K89").Value = 8100.666499999999
$ws.Range("L89").Value = 15041.6665
$ws.Range("M89").Value = -2484.666499999999
$ws.Range("N89").Value = -26273.6665
$ws.Range("H134").Value = 43723.5
$ws.Range("I134").Value = 2226.0952
$ws.Range("J134").Value = 334205.34
$ws.Range("K134").Value = 6678.285600000001
$ws.Range("L134").Value = 1002616.02
$ws.Range("M134").Value = -4143.285600000001
$ws.Range("N134").Value = -1007686.02

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 487044.03
$ws.Range("I31").Value = 7946.2964
$ws.Range("J31").Value = 762270.4
$ws.Range("K31").Value = 7946.2964
$ws.Range("L31").Value = 762270.4
$ws.Range("M31").Value = -7651.2964
$ws.Range("N31").Value = -762860.4
$ws.Range("H34").Value = 487044.03
$ws.Range("I34").Value = 7946.2964
$ws.Range("J34").Value = 762270.4
$ws.Range("K34").Value = 7946.2964
$ws.Range("L34").Value = 762270.4
$ws.Range("M34").Value = -7744.2964
$ws.Range("N34").Value = -762674.4
$ws.Range("H58").Value = 3035977.5
$ws.Range("I58").Value = 4547443.5
$ws.Range("J58").Value = 13045.3
$ws.Range("K58").Value = 4547443.5
$ws.Range("L58").Value = 13045.3
$ws.Range("M58").Value = -4547240.5
$ws.Range("N58").Value = -13451.3
$ws.Range("H122").Value = 4710.273
$ws.Range("I122").Value = 4285.5713
$ws.Range("J122").Value = 5453.5
$ws.Range("K122").Value = 12856.7139
$ws.Range("L122").Value = 16360.5
$ws.Range("M122").Value = -10406.7139
$ws.Range("N122").Value = -21260.5
$ws.Range("H136").Value = 3035977.5
$ws.Range("I136").Value = 4547443.5
$ws.Range("J136").Value = 13045.3
$ws.Range("K136").Value = 13642330.5
$ws.Range("L136").Value = 39135.89999999999
$ws.Range("M136").Value = -13639780.5
$ws.Range("N136").Value = -44235.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 208.33333
$ws.Range("I11").Value = 159.41667
$ws.Range("J11").Value = 404
$ws.Range("K11").Value = 478.25001
$ws.Range("L11").Value = 1212
$ws.Range("M11").Value = -338.25001
$ws.Range("N11").Value = -1492
$ws.Range("H68").Value = 118550.09
$ws.Range("I68").Value = 240518.23
$ws.Range("J68").Value = 3900.04
$ws.Range("K68").Value = 721554.6900000001
$ws.Range("L68").Value = 11700.12
$ws.Range("M68").Value = -720743.6900000001
$ws.Range("N68").Value = -13322.12
$ws.Range("H71").Value = 118550.09
$ws.Range("I71").Value = 240518.23
$ws.Range("J71").Value = 3900.04
$ws.Range("K71").Value = 2164664.07
$ws.Range("L71").Value = 35100.36
$ws.Range("M71").Value = -2160608.07
$ws.Range("N71").Value = -43212.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5822.615
$ws.Range("I122").Value = 7650.75
$ws.Range("K122").Value = 22952.25
$ws.Range("M122").Value = -20502.25
$ws.Range("H132").Value = 12487.385
$ws.Range("I132").Value = 5506
$ws.Range("K132").Value = 16518
$ws.Range("M132").Value = -13988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H132").Value = 21256.334
$ws.Range("I132").Value = 26044
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 78132
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -75602
$ws.Range("N132").Value = -18558.5
$ws.Range("H136").Value = 3760.6
$ws.Range("I136").Value = 2366.5386
$ws.Range("J136").Value = 5270.8335
$ws.Range("K136").Value = 7099.6158
$ws.Range("L136").Value = 15812.5005
$ws.Range("M136").Value = -4549.6158
$ws.Range("N136").Value = -20912.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 28575358
$ws.Range("I81").Value = 3749
$ws.Range("J81").Value = 40004000
$ws.Range("K81").Value = 7498
$ws.Range("L81").Value = 80008000
$ws.Range("M81").Value = -6437
$ws.Range("N81").Value = -80010122
$ws.Range("H84").Value = 28575358
$ws.Range("I84").Value = 3749
$ws.Range("J84").Value = 40004000
$ws.Range("K84").Value = 37490
$ws.Range("L84").Value = 400040000
$ws.Range("M84").Value = -32186
$ws.Range("N84").Value = -400050608
$ws.Range("H132").Value = 3368.7727
$ws.Range("I132").Value = 3425.625
$ws.Range("K132").Value = 10276.875
$ws.Range("M132").Value = -7746.875
$ws.Range("H136").Value = 6603.705
$ws.Range("I136").Value = 7094.9644
$ws.Range("K136").Value = 21284.8932
$ws.Range("M136").Value = -18734.8932
